# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
# These reflect a re-run of the scrape at a later time (03:0x -> 07:0x UTC).
$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshot.Range("K2").Value = "2025-12-20T07:04:46.731995+00:00"
$snapshot.Range("K3").Value = "2025-12-20T07:04:49.474225+00:00"
$snapshot.Range("K4").Value = "2025-12-20T07:04:49.474259+00:00"
$snapshot.Range("K5").Value = "2025-12-20T07:04:49.474280+00:00"
$snapshot.Range("K6").Value = "2025-12-20T07:04:51.674165+00:00"
$snapshot.Range("K7").Value = "2025-12-20T07:04:54.391848+00:00"
$snapshot.Range("K8").Value = "2025-12-20T07:04:56.702164+00:00"
$snapshot.Range("K9").Value = "2025-12-20T07:04:56.702190+00:00"
$snapshot.Range("K10").Value = "2025-12-20T07:04:58.891352+00:00"
$snapshot.Range("K11").Value = "2025-12-20T07:05:03.963923+00:00"
$snapshot.Range("K12").Value = "2025-12-20T07:05:06.603479+00:00"
$snapshot.Range("K13").Value = "2025-12-20T07:05:09.360728+00:00"
$snapshot.Range("K14").Value = "2025-12-20T07:05:11.646777+00:00"
$snapshot.Range("K15").Value = "2025-12-20T07:05:13.925021+00:00"
$snapshot.Range("K16").Value = "2025-12-20T07:05:13.925050+00:00"
$snapshot.Range("K17").Value = "2025-12-20T07:05:13.925067+00:00"
$snapshot.Range("K18").Value = "2025-12-20T07:05:13.925085+00:00"
$snapshot.Range("K19").Value = "2025-12-20T07:05:16.181666+00:00"
$snapshot.Range("K20").Value = "2025-12-20T07:05:16.181700+00:00"
$snapshot.Range("K21").Value = "2025-12-20T07:05:16.181725+00:00"
$snapshot.Range("K22").Value = "2025-12-20T07:05:29.568073+00:00"
$snapshot.Range("K23").Value = "2025-12-20T07:05:29.568103+00:00"
$snapshot.Range("K24").Value = "2025-12-20T07:05:29.568123+00:00"
$snapshot.Range("K25").Value = "2025-12-20T07:05:32.300445+00:00"
$snapshot.Range("K26").Value = "2025-12-20T07:05:32.300471+00:00"
$snapshot.Range("K27").Value = "2025-12-20T07:05:32.300488+00:00"
$snapshot.Range("K28").Value = "2025-12-20T07:05:34.583579+00:00"
$snapshot.Range("K29").Value = "2025-12-20T07:05:34.583608+00:00"
$snapshot.Range("K30").Value = "2025-12-20T07:05:34.583627+00:00"
$snapshot.Range("K31").Value = "2025-12-20T07:05:36.793663+00:00"
$snapshot.Range("K32").Value = "2025-12-20T07:05:39.125522+00:00"
$snapshot.Range("K33").Value = "2025-12-20T07:05:39.125554+00:00"
$snapshot.Range("K34").Value = "2025-12-20T07:05:43.582991+00:00"
$snapshot.Range("K35").Value = "2025-12-20T07:05:43.583021+00:00"
$snapshot.Range("K36").Value = "2025-12-20T07:05:45.856236+00:00"
$snapshot.Range("K37").Value = "2025-12-20T07:05:45.856270+00:00"

# Remove the now-stale "new_injured" entry: the player (Samsonov Ilya /
# СОЧ) who was newly injured in the earlier run is no longer reported as
# a fresh injury in the later run, so the single data row is deleted from
# the "new_injured" sheet. This shifts the sheet dimension from A1:G2 to
# A1:G1 (header row only remains).
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
